$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Severity" column (D) with RPN-ish values, header bold ---
$ws.Range("D1").Value = "Severity"
$ws.Range("D2").Value = 22
$ws.Range("D3").Value = 42
$ws.Range("D4").Value = 9
$ws.Range("D5").Value = 91
$ws.Range("D6").Value = 61
$ws.Range("D7").Value = 44

# --- Reformat C2:C7 (integer thousand-separator format, right aligned, normal font) ---
$ws.Range("C2:C7").NumberFormat = "#,##0"
$ws.Range("C2:C7").HorizontalAlignment = -4152
$ws.Range("C2:C7").Font.Name = "Calibri"
$ws.Range("C2:C7").Font.Size = 11

# --- Clear the old mirrored values that used to live in column E, giving
#     the (now empty) cells the exact same formatting as C2:C7 ---
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E2:E4").NumberFormat = "#,##0"
$ws.Range("E2:E4").HorizontalAlignment = -4152
$ws.Range("E2:E4").Font.Name = "Calibri"
$ws.Range("E2:E4").Font.Size = 11

# --- Header D1 bold ---
$ws.Range("D1").Font.Bold = $true

# --- Row heights: header keeps 19.5, data rows now 18.75 ---
$ws.Rows(2).RowHeight = 18.75
$ws.Rows(3).RowHeight = 18.75
$ws.Rows(4).RowHeight = 18.75
$ws.Rows(5).RowHeight = 18.75
$ws.Rows(6).RowHeight = 18.75
$ws.Rows(7).RowHeight = 18.75

# --- New rows 9-10, underlined empty cells in D (RPN notes area) ---
$ws.Range("D9").Font.Underline = $true
$ws.Range("D10").Font.Underline = $true

# --- Selection as left by the author ---
$ws.Range("D9:D10").Select()
